# "render termino de version entregable"
# Rename the sheet, add header/data styling (fills, fonts, borders, alignment,
# number formats), freeze the header row, hide gridlines, set column widths
# and turn on an AutoFilter over the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet --------------------------------------------------------
$ws.Name = "Responsable"

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth  = 11.71
$ws.Columns.Item(2).ColumnWidth  = 7.71
$ws.Columns.Item(3).ColumnWidth  = 13.71
$ws.Columns.Item(4).ColumnWidth  = 15.71
$ws.Columns.Item(5).ColumnWidth  = 17.71
$ws.Columns.Item(6).ColumnWidth  = 19.71
$ws.Columns.Item(7).ColumnWidth  = 9.71
$ws.Columns.Item(8).ColumnWidth  = 10.71
$ws.Columns.Item(9).ColumnWidth  = 12.71
$ws.Columns.Item(10).ColumnWidth = 22.71
$ws.Columns.Item(11).ColumnWidth = 8.71
$ws.Columns.Item(12).ColumnWidth = 8.71
$ws.Columns.Item(13).ColumnWidth = 8.71
$ws.Columns.Item(14).ColumnWidth = 8.71
$ws.Columns.Item(15).ColumnWidth = 8.71
$ws.Columns.Item(16).ColumnWidth = 8.71

# --- Header row styling (row 1): bold white text on dark-blue fill -------
$header = $ws.Range("A1:P1")
$header.Font.Name = "Calibri"
$header.Font.Bold = $true
$header.Font.Color = 16777215          # white
$header.Interior.Color = 7949855       # 1F4E79 (dark blue)
$header.HorizontalAlignment = -4108    # xlCenter
$header.VerticalAlignment = -4108      # xlCenter
$header.Borders.Item(9).LineStyle = 1  # xlEdgeBottom
$header.Borders.Item(9).Weight = 2     # xlThin
$header.Borders.Item(9).Color = 0

# --- Data rows styling (rows 2-9): black text, thin box border ----------
$data = $ws.Range("A2:P9")
$data.Font.Name = "Calibri"
$data.Font.Bold = $false
$data.Font.Color = 0                   # black
$data.VerticalAlignment = -4108        # xlCenter
$data.Borders.LineStyle = 1
$data.Borders.Weight = 2
$data.Borders.Color = 0

# column A (names) -> left aligned
$ws.Range("A2:A9").HorizontalAlignment = -4131   # xlLeft
# counts columns -> centered, #,##0
$ws.Range("B2:C9").HorizontalAlignment = -4108
$ws.Range("B2:C9").NumberFormat = "#,##0"
$ws.Range("E2:E9,G2:H9,J2:J9").HorizontalAlignment = -4108
$ws.Range("E2:E9,G2:H9,J2:J9").NumberFormat = "#,##0"
# percent columns -> centered, 0.0"%"
$ws.Range("D2:D9,F2:F9,I2:I9").HorizontalAlignment = -4108
$ws.Range("D2:D9,F2:F9,I2:I9").NumberFormat = "0.0""%"""
# remaining stat columns (K:P) -> centered, default number format
$ws.Range("K2:P9").HorizontalAlignment = -4108

# --- View: hide gridlines, freeze header row -----------------------------
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- AutoFilter over the used range --------------------------------------
$used = $ws.Range("A1:P9")
$used.AutoFilter() | Out-Null
$fname = $ws.Names.Add("_xlnm._FilterDatabase", "='Responsable'!`$A`$1:`$P`$9")
$fname.Visible = $false

$ws.Range("A1").Select()
